$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear J2 (was "Number of nodes"), keep formatting/style
$ws.Range("J2").ClearContents()

# Remove values in C18:C20 (5000, 10000, 20000)
$ws.Range("C18:C20").ClearContents()

# Update the selection to J2 (matches sheetView selection change)
$ws.Range("J2").Select()
